$d = $word.ActiveDocument

# Find.Execute positional signature:
#   (FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#    MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# wdFindContinue = 1, wdReplaceAll = 2

# 1) Update the month label from "MONTH: TOTAL" to "MONTH: JUNE-OCT"
$d.Content.Find.Execute("TOTAL", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "JUNE-OCT", 2)

# 2) In the summary table, change every lone "+0" gap cell to "0".
#    MatchWholeWord keeps the "20 + 20" / "20 + 0" total cells untouched.
$d.Content.Find.Execute("+0", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "0", 2)
